$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "Акция" best-offer badge, header styled like the other
# header cells (bold, centered - matches the "Картинка" header look).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Акция"

# Body cells F2:F10 pick up the same (non-bold) font as the rest of the
# data rows.
$ws.Range("E2").Copy()
$ws.Range("F2:F10").PasteSpecial(-4122)
$ws.Range("F2:F10").HorizontalAlignment = 1

# Only rows with a defined best offer show the badge text; the rest of
# the column stays blank.
$ws.Range("F2").Value = "Выгодное предложение"
$ws.Range("F9").Value = "Выгодное предложение"

$ws.Range("F1:F10").EntireColumn.AutoFit()

$excel.CutCopyMode = $false
